$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Resmi_Tatiller")
$ws2.Select()
$ws1 = $wb.Worksheets.Item("Belirli_Gun_ve_Haftalar")
$ws1.Select()
